$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Branch LPF")

# Update Gain value (B18): 4 -> 2
$ws.Range("B18").Value = 2

# Set the active selection to D20 as captured in the saved view state
$ws.Activate()
$ws.Range("D20").Select()

$wb.Save()
